$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-9 from 2023-10-25 (45224)
# to 2023-11-03 (45233), keeping the existing date formatting.
$ws.Range("C2:C9").Value = 45233
